$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns C (price), E (reviews) and F (sku) can look numeric to Excel's
# smart-parsing ($629.99, 6615769, "1,275", ...), which would silently turn
# them into Number cells. Mark the ranges as Text first so the assigned
# values stick as plain text (matching the source t="inlineStr" cells);
# the temporary formatting is cleared again at the end so no stray number
# format/style is left behind on the cells.
$txtC = $ws.Range("C2:C15")
$txtC21 = $ws.Range("C21")
$txtE = $ws.Range("E2:E21")
$txtF = $ws.Range("F2:F15")
$txtF21 = $ws.Range("F21")
$txtC.NumberFormat = "@"
$txtC21.NumberFormat = "@"
$txtE.NumberFormat = "@"
$txtF.NumberFormat = "@"
$txtF21.NumberFormat = "@"

# --- Row-by-row refresh of the scraped product listing (A:G) ---

$ws.Range("A2").Value = 'Lenovo-  Yoga 7i 2 - in - 1  -  Copilot+ PC  -  16" 2K Touchscreen Laptop  -  Intel Core Ultra 5 Processor  -  16GB Memory  -  512GB SSD  -  Luna Grey'
$ws.Range("B2").Value = 'https://www.bestbuy.com/site/lenovo-yoga-7i-2-in-1-copilot-pc-16-2k-touchscreen-laptop-intel-core-ultra-5-processor-16gb-memory-512gb-ssd-luna-grey/6615769.p?skuId=6615769'
$ws.Range("C2").Value = '$649.99'
$ws.Range("D2").Value = 'Rating 4.9 out of 5 stars with 7 reviews'
$ws.Range("E2").Value = '7'
$ws.Range("F2").Value = '6615769'
$ws.Range("G2").Value = '83JT0000US'

$ws.Range("A3").Value = 'HP-  OmniBook 5 Flip 2 - in - 1 14" 2K Touch - Screen Laptop  -  Intel Core 7  -  16GB Memory  -  512GB SSD  -  Glacier Silver'
$ws.Range("B3").Value = 'https://www.bestbuy.com/site/hp-omnibook-5-flip-2-in-1-14-2k-touch-screen-laptop-intel-core-7-16gb-memory-512gb-ssd-glacier-silver/6614107.p?skuId=6614107'
$ws.Range("C3").Value = '$549.99'
$ws.Range("D3").Value = 'Rating 4.8 out of 5 stars with 19 reviews'
$ws.Range("E3").Value = '19'
$ws.Range("F3").Value = '6614107'
$ws.Range("G3").Value = '14-fp0023dx'

$ws.Range("A4").Value = 'Lenovo-  Yoga 9i 2 - in - 1 14" 2.8K OLED Touchscreen Laptop with Pen  -  Intel Core Ultra 7 155H with 16GB Memory  -  1TB SSD  -  Cosmic Blue'
$ws.Range("B4").Value = 'https://www.bestbuy.com/site/lenovo-yoga-9i-2-in-1-14-2-8k-oled-touchscreen-laptop-with-pen-intel-core-ultra-7-155h-with-16gb-memory-1tb-ssd-cosmic-blue/6571371.p?skuId=6571371'
$ws.Range("C4").Value = '$999.99'
$ws.Range("D4").Value = 'Rating 4.4 out of 5 stars with 222 reviews'
$ws.Range("E4").Value = '222'
$ws.Range("F4").Value = '6571371'
$ws.Range("G4").Value = '83AC0001US'

$ws.Range("A5").Value = 'HP-  17.3" Full HD Laptop  -  AMD Ryzen 5  -  8GB Memory  -  512GB SSD  -  Natural Silver'
$ws.Range("B5").Value = 'https://www.bestbuy.com/site/hp-17-3-full-hd-laptop-amd-ryzen-5-8gb-memory-512gb-ssd-natural-silver/6612252.p?skuId=6612252'
$ws.Range("C5").Value = '$629.99'
$ws.Range("D5").Value = 'Rating 4.7 out of 5 stars with 216 reviews'
$ws.Range("E5").Value = '216'
$ws.Range("F5").Value = '6612252'
$ws.Range("G5").Value = '17-cp2025dx'

$ws.Range("A6").Value = 'HP-  Victus 15.6" 144Hz Full HD Gaming Laptop  -  Intel Core i5  -  8GB Memory  -  NVIDIA GeForce RTX 3050  -  512GB SSD  -  Mica Silver'
$ws.Range("B6").Value = 'https://www.bestbuy.com/site/hp-victus-15-6-144hz-full-hd-gaming-laptop-intel-core-i5-8gb-memory-nvidia-geforce-rtx-3050-512gb-ssd-mica-silver/6618924.p?skuId=6618924'
$ws.Range("C6").Value = '$529.99'
$ws.Range("D6").Value = 'Rating 4.5 out of 5 stars with 50 reviews'
$ws.Range("E6").Value = '50'
$ws.Range("F6").Value = '6618924'
$ws.Range("G6").Value = '15-fa2013dx'

$ws.Range("A7").Value = 'Lenovo-  Yoga 7i 2 - in - 1  -  Copilot+ PC  -  14" 2K OLED Touchscreen Laptop  -  Intel Core Ultra 7 Processor  -  16GB Memory  -  1TB SSD  -  Luna Grey'
$ws.Range("B7").Value = 'https://www.bestbuy.com/site/lenovo-yoga-7i-2-in-1-copilot-pc-14-2k-oled-touchscreen-laptop-intel-core-ultra-7-processor-16gb-memory-1tb-ssd-luna-grey/6615777.p?skuId=6615777'
$ws.Range("C7").Value = '$799.99'
$ws.Range("D7").Value = 'Rating 5 out of 5 stars with 1 reviewfalse'
$ws.Range("E7").Value = '1'
$ws.Range("F7").Value = '6615777'
$ws.Range("G7").Value = '83JQ000KUS'

$ws.Range("A8").Value = 'HP-  14" Refurbished 1920 x 1080 FHD  -  Intel 11th Gen Core i5 - 1145G7 with 32GB RAM  -  Intel Iris Xe Graphics  -  1TB SSD  -  Silver'
$ws.Range("B8").Value = 'https://www.bestbuy.com/site/hp-14-refurbished-1920-x-1080-fhd-intel-11th-gen-core-i5-1145g7-with-32gb-ram-intel-iris-xe-graphics-1tb-ssd-silver/6545476.p?skuId=6545476'
$ws.Range("C8").Value = '$666.99'
$ws.Range("D8").Value = 'Rating 5 out of 5 stars with 1 reviewfalse'
$ws.Range("E8").Value = '1'
$ws.Range("F8").Value = '6545476'
$ws.Range("G8").Value = '840 G8'

$ws.Range("A9").Value = 'HP-  ZBook 15U G5 15.6" Refurbished Laptop  -  Intel 8th Gen Core i7 with 32GB Memory  -  Intel UHD Graphics 620  -  512GB SSD  -  Silver'
$ws.Range("B9").Value = 'https://www.bestbuy.com/site/hp-zbook-15u-g5-15-6-refurbished-laptop-intel-8th-gen-core-i7-with-32gb-memory-intel-uhd-graphics-620-512gb-ssd-silver/6579728.p?skuId=6579728'
$ws.Range("C9").Value = '$643.99'
$ws.Range("D9").Value = 'Not yet reviewed'
$ws.Range("E9").Value = 'Not yet reviewed'
$ws.Range("F9").Value = '6579728'
$ws.Range("G9").Value = 'ZBook 15U G5'

$ws.Range("A10").Value = 'HP-  ProBook 440 G11 14" IPS 1920 x 1200 (WUXGA) Laptop  -  Intel Core Ultra 5 with 16GB Memory  -  256 GB SSD  -  Pike Silver, Silver'
$ws.Range("B10").Value = 'https://www.bestbuy.com/site/hp-probook-440-g11-14-ips-1920-x-1200-wuxga-laptop-intel-core-ultra-5-with-16gb-memory-256-gb-ssd-pike-silver-silver/6588386.p?skuId=6588386'
$ws.Range("C10").Value = '$1,089.99'
$ws.Range("D10").Value = 'Not yet reviewed'
$ws.Range("E10").Value = 'Not yet reviewed'
$ws.Range("F10").Value = '6588386'
$ws.Range("G10").Value = 'A1LC2UT#ABA'

$ws.Range("A11").Value = 'HP-  Envy 2 - in - 1 16" 2K Touch - Screen Laptop  -  AMD Ryzen 7  -  16GB Memory  -  1TB SSD  -  Meteor Silver'
$ws.Range("B11").Value = 'https://www.bestbuy.com/site/hp-envy-2-in-1-16-2k-touch-screen-laptop-amd-ryzen-7-16gb-memory-1tb-ssd-meteor-silver/6571083.p?skuId=6571083'
$ws.Range("C11").Value = '$760.99'
$ws.Range("D11").Value = 'Rating 4.7 out of 5 stars with 790 reviews'
$ws.Range("E11").Value = '790'
$ws.Range("F11").Value = '6571083'
$ws.Range("G11").Value = '16-ad0023dx/9S1M4UA#ABA'

$ws.Range("A12").Value = 'HP-  OmniBook Ultra Flip  -  Copilot+ PC  -  14" 3K OLED Touch - Screen Laptop  -  Intel Core Ultra 7  -  16GB Memory  -  1TB SSD  -  Eclipse Grey'
$ws.Range("B12").Value = 'https://www.bestbuy.com/site/hp-omnibook-ultra-flip-copilot-pc-14-3k-oled-touch-screen-laptop-intel-core-ultra-7-16gb-memory-1tb-ssd-eclipse-grey/6593552.p?skuId=6593552'
$ws.Range("C12").Value = '$1,199.99'
$ws.Range("D12").Value = 'Rating 4.5 out of 5 stars with 155 reviews'
$ws.Range("E12").Value = '155'
$ws.Range("F12").Value = '6593552'
$ws.Range("G12").Value = '14-fh0013dx/A9SR3UA#ABA'

$ws.Range("A13").Value = 'HP-  Envy 2 - in - 1 16" 2K Touch - Screen Laptop  -  Intel Core Ultra 7  -  16GB Memory  -  1TB SSD  -  Glacier Silver'
$ws.Range("B13").Value = 'https://www.bestbuy.com/site/hp-envy-2-in-1-16-2k-touch-screen-laptop-intel-core-ultra-7-16gb-memory-1tb-ssd-glacier-silver/6571084.p?skuId=6571084'
$ws.Range("C13").Value = '$760.99'
$ws.Range("D13").Value = 'Rating 4.7 out of 5 stars with 1275 reviews'
$ws.Range("E13").Value = '1,275'
$ws.Range("F13").Value = '6571084'
$ws.Range("G13").Value = '16-ac0023dx/9S1R6UA#ABA'

$ws.Range("A14").Value = 'HP-  Envy 2 - in - 1 14" 2K Touch - Screen Laptop  -  Intel Core Ultra 7  -  16GB Memory  -  1TB SSD  -  Meteor Silver'
$ws.Range("B14").Value = 'https://www.bestbuy.com/site/hp-envy-2-in-1-14-2k-touch-screen-laptop-intel-core-ultra-7-16gb-memory-1tb-ssd-meteor-silver/6571085.p?skuId=6571085'
$ws.Range("C14").Value = '$701.99'
$ws.Range("D14").Value = 'Rating 4.7 out of 5 stars with 569 reviews'
$ws.Range("E14").Value = '569'
$ws.Range("F14").Value = '6571085'
$ws.Range("G14").Value = '14-fc0023dx/9T8G4UA#ABA'

$ws.Range("A15").Value = 'HP-  Envy 2 - in - 1 16" 2K Touch - Screen Laptop  -  Intel Core Ultra 5  -  16GB Memory  -  512GB SSD  -  Glacier Silver'
$ws.Range("B15").Value = 'https://www.bestbuy.com/site/hp-envy-2-in-1-16-2k-touch-screen-laptop-intel-core-ultra-5-16gb-memory-512gb-ssd-glacier-silver/6571081.p?skuId=6571081'
$ws.Range("C15").Value = '$617.99'
$ws.Range("D15").Value = 'Rating 4.7 out of 5 stars with 525 reviews'
$ws.Range("E15").Value = '525'
$ws.Range("F15").Value = '6571081'
$ws.Range("G15").Value = '16-ac0013dx/9S1R5UA#ABA'

$ws.Range("A21").Value = 'Lenovo-  Yoga 7i 2 - in - 1 16" 2K Touchscreen Laptop  -  Intel Core Ultra 7 155U with 16GB Memory  -  1TB SSD  -  Storm Grey'
$ws.Range("B21").Value = 'https://www.bestbuy.com/site/lenovo-yoga-7i-2-in-1-16-2k-touchscreen-laptop-intel-core-ultra-7-155u-with-16gb-memory-1tb-ssd-storm-grey/6571369.p?skuId=6571369'
$ws.Range("C21").Value = '$1,049.99'
$ws.Range("D21").Value = 'Rating 4.7 out of 5 stars with 1609 reviews'
$ws.Range("E21").Value = '1,609'
$ws.Range("F21").Value = '6571369'
$ws.Range("G21").Value = '83DL0002US'

# Rows 16-20 keep their "N/A" placeholders, but the reviews count (E) now
# round-trips as the literal text "0" instead of the numeric 0.
$ws.Range("E16").Value = "0"
$ws.Range("E17").Value = "0"
$ws.Range("E18").Value = "0"
$ws.Range("E19").Value = "0"
$ws.Range("E20").Value = "0"

# Drop the temporary text-forcing format now that every value is written,
# so the cells end up unstyled again (same as the source file).
$txtC.ClearFormats()
$txtC21.ClearFormats()
$txtE.ClearFormats()
$txtF.ClearFormats()
$txtF21.ClearFormats()

# The source now only has 20 data rows (2-21); drop the trailing all-N/A
# rows 22-24 so the sheet dimension shrinks from A1:G24 to A1:G21.
$ws.Rows("22:24").Delete()
